$wb = $excel.ActiveWorkbook

$wsAddr = $wb.Worksheets.Item("Station addresses")
$wsLink = $wb.Worksheets.Item("Linking station names")

# "Station addresses": remove the throwaway "Stuttgart" test row (was row 74)
$wsAddr.Rows.Item(74).Delete()

# "Linking station names": remove the matching lookup row that pointed at it (was row 81)
$wsLink.Rows.Item(81).Delete()

# Leave the selection/active sheet the way the author's session ended up:
# focus back on the row that used to hold "Stuttgart" in the address sheet,
# and park the linking sheet's selection on D78.
$wsLink.Activate() | Out-Null
$wsLink.Range("D78").Select() | Out-Null

$wsAddr.Activate() | Out-Null
$wsAddr.Rows.Item(74).Select() | Out-Null
